$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.444.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "'1.817.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.21%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'318.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").Value = "'0.9981"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.5727"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.62%  "
$ws.Range("D8").Value = "'0.3854"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.88%  "
$ws.Range("D9").Value = "'43.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.07646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("E11").Value = "  +8.37%  "
$ws.Range("D12").Value = "'21.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.37%  "
$ws.Range("D13").Value = "'0.9977"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'6.261"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("D15").Value = "'1.807.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.70%  "
$ws.Range("D16").Value = "'7.284"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.74%  "
$ws.Range("D17").Value = "'92.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.21%  "
$ws.Range("E18").Value = "  +4.83%  "
$ws.Range("D19").Value = "'0.06519"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "'0.9979"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "'17.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.66%  "
$ws.Range("D22").Value = "'6.014"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.14%  "
$ws.Range("D23").Value = "'28.459.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").Value = "'2.092"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'20.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.13%  "
$ws.Range("D27").Value = "'157.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  +16.13%  "
$ws.Range("D29").Value = "'2.019.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").Value = "'123.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").Value = "'1.157"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.43%  "
$ws.Range("D32").Value = "'0.1055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.33%  "
$ws.Range("D33").Value = "'5.798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.75%  "
$ws.Range("D34").Value = "'3.638"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "'0.02323"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("D36").Value = "'0.2169"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.40%  "
$ws.Range("D37").Value = "'8.776"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.46%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.6511"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.78%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("D41").Value = "'0.06102"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").Value = "'0.9977"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'1.157"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("D44").Value = "'1.379"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "'13.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.51%  "
$ws.Range("D46").Value = "'0.6048"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.61%  "
$ws.Range("D47").Value = "'3.721"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("D48").Value = "'122.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'1.952"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("D50").Value = "'1.150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").Value = "'0.06856"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.91%  "
